# Add "verlenging_bij_vernieuwing" concept to the codelijst (LVBR-379)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 6 (the "vernieuwing" row), shifting it and all
# following rows down by one.
$ws.Rows(6).Insert()

# Populate the newly inserted row 6 with the new concept's data.
$ws.Range("A6").Value = "https://data.omgeving.vlaanderen.be/id/concept/leegstand/aanvraag_opschorting_heffing_reden/verlenging_bij_vernieuwing"
$ws.Range("B6").Value = "http://www.w3.org/2004/02/skos/core#Concept"
$ws.Range("C6").Value = "https://data.omgeving.vlaanderen.be/id/concept/leegstand/procedure/aanvraag_opschorting_heffing"
$ws.Range("D6").Value = "https://data.omgeving.vlaanderen.be/id/conceptscheme/leegstand/aanvraag_opschorting_heffing_reden"
$ws.Range("E6").Value = "verlenging_bij_vernieuwing"
$ws.Range("F6").Value = "Verlenging bij vernieuwing"
$ws.Range("G6").Value = "https://data.omgeving.vlaanderen.be/id/conceptscheme/leegstand/aanvraag_opschorting_heffing_reden"
$ws.Range("H6").Value = "null"
$ws.Range("I6").Value = "null"
$ws.Range("J6").Value = "null"
$ws.Range("K6").Value = "null"
$ws.Range("L6").Value = "null"
$ws.Range("M6").Value = "null"
$ws.Range("N6").Value = "null"
$ws.Range("O6").Value = "null"
$ws.Range("P6").Value = "null"
$ws.Range("Q6").Value = "null"
$ws.Range("R6").Value = "null"
$ws.Range("S6").Value = "null"

# The conceptscheme row for "aanvraag_opschorting_heffing_reden" (previously row
# 40) is now row 41 after the insert. Update its hasTopConcept (S) list to
# include the newly added concept, in alphabetical order.
$ws.Range("S41").Value = "https://data.omgeving.vlaanderen.be/id/concept/leegstand/aanvraag_opschorting_heffing_reden/bodemsaneringsproject|https://data.omgeving.vlaanderen.be/id/concept/leegstand/aanvraag_opschorting_heffing_reden/brownfieldconvenant|https://data.omgeving.vlaanderen.be/id/concept/leegstand/aanvraag_opschorting_heffing_reden/leegstaand_niet_verwaarloosd|https://data.omgeving.vlaanderen.be/id/concept/leegstand/aanvraag_opschorting_heffing_reden/nieuwe_eigenaars|https://data.omgeving.vlaanderen.be/id/concept/leegstand/aanvraag_opschorting_heffing_reden/verlenging_bij_vernieuwing|https://data.omgeving.vlaanderen.be/id/concept/leegstand/aanvraag_opschorting_heffing_reden/vernieuwing"
